$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 469, shifting existing rows 469:488 down to 470:489
$ws.Rows.Item(469).Insert()

# Populate the newly inserted row 469 with the new weekly data point
$ws.Range("A469").Value = 3
$ws.Range("B469").Value = "Femacal de La Calera"
$ws.Range("C469").Value = "Coquimbo"
$ws.Range("D469").Value = 44939
$ws.Range("E469").Value = 5
$ws.Range("F469").Value = 100112040
$ws.Range("G469").Value = "Cilantro"
$ws.Range("H469").Value = "Sin especificar"
$ws.Range("I469").Value = "Primera"
$ws.Range("J469").Value = 120
$ws.Range("K469").Value = 4000
$ws.Range("L469").Value = 4000
$ws.Range("M469").Value = 4000
$ws.Range("N469").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O469").Value = "Provincia de Quillota"
$ws.Range("P469").Value = 1333
$ws.Range("Q469").Value = 3
$ws.Range("R469").Value = "Hortaliza"
